$wb = $excel.ActiveWorkbook
$conv = $wb.Worksheets.Item("conversion")
$cdiac = $wb.Worksheets.Item("CDIAC")

# Fix fraction oxidized for liquid fuel: was "=0.918" formula, now plain constant 0.985
$conv.Range("B21").Value = 0.985

# CDIAC sheet: tie B3:B6 to the conversion sheet fraction-oxidized values via formulas
$cdiac.Range("B3").Formula = "=0.85*conversion!B21"
$cdiac.Range("B4").Formula = "=0.0137*conversion!B22"
$cdiac.Range("B5").Formula = "=0.855*conversion!B23"
$cdiac.Range("B6").Formula = "=13.454*conversion!B24"

# Fix source citation cells G3/G4 to match G5/G6 (ndp030 table4 URL), removing special styling
$cdiac.Range("G3").Value = "http://cdiac.ornl.gov/epubs/ndp/ndp030/tables/table4.htm"
$cdiac.Range("G4").Value = "http://cdiac.ornl.gov/epubs/ndp/ndp030/tables/table4.htm"
$cdiac.Range("G3").ClearFormats()
$cdiac.Range("G4").ClearFormats()

# Shorten the fraction-oxidized source citation on the conversion sheet
$conv.Range("A19").Value = "Source: http://cdiac.ornl.gov/epubs/ndp/ndp030/tables/table4.htm"
